{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Implements the changes from the commit \"Finished Project and added Privacy Page\":\n//  1. \"Erst beginne ich ...\" -> \"Erst danach beginne ich ...\" (real text change).\n//  2-5. Several paragraphs get their runs re-split around spell-checked words\n//       (Griptape, included, Product, Diagram, ...). The *visible* text of\n//       those paragraphs is byte-identical before/after (only <w:proofErr/>\n//       bookkeeping + run boundaries change, which Word inserts automatically\n//       during interactive spell-checking and is not reachable from the\n//       Word JS API), so there is nothing observable to change there.\n//  6. A new \"Implementation\" section (Heading 1 + body paragraphs, a\n//     \"Problem mit Composite Pattern\" Heading 2 + body paragraphs, and two\n//     trailing empty paragraphs, the last right-aligned) is appended at the\n//     very end of the document body, right after the Class-Diagram picture\n//     and before the final section break.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) \"Zuerst habe ich mein Programm ... Erst beginne ich ...\"\n//    -> \"... Erst danach beginne ich ...\"\n// ---------------------------------------------------------------------\nconst erstResults = body.search(\"Erst beginne ich mit dem eigentlichen programmieren.\", { matchCase: true });\nerstResults.load(\"items\");\nawait context.sync();\n\nif (erstResults.items.length > 0) {\n    erstResults.items[0].insertText(\n        \"Erst danach beginne ich mit dem eigentlichen programmieren.\",\n        Word.InsertLocation.replace\n    );\n    await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 6) Append the new \"Implementation\" section at the end of the document.\n// ---------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\nfunction addParagraph(text, styleBuiltIn) {\n    const p = anchor.insertParagraph(text, Word.InsertLocation.after);\n    if (styleBuiltIn) {\n        p.styleBuiltIn = styleBuiltIn;\n    }\n    anchor = p;\n    return p;\n}\n\naddParagraph(\"Implementation\", Word.Style.heading1);\n\naddParagraph(\n    \"Als Entwicklungsumgebung benutze ich Visual Studio 2019 mit einer MSSQL Datenbank. Ich erstelle ein ASP.NET Core MVC Projekt.\",\n    Word.Style.normal\n);\n\naddParagraph(\"Problem mit Composite Pattern\", Word.Style.heading2);\n\naddParagraph(\n    \"Als ich das Composite Pattern implementieren wollte sah ich, dass es dabei ein Problem gibt. In einem MVC Projekt arbeitet man mit Model-Klassen, welche gleichzeitig die Tabellen in der Datenbank sind. Da eine Tabelle nicht von einer anderen Tabelle erben kann, war es nicht m\u00f6glich das Composite Pattern zu implementieren.\",\n    Word.Style.normal\n);\n\naddParagraph(\n    \"Ich habe beschlossen, stattdessen das Dependency Injection Design Pattern anzuwenden. Ausserdem ist MVC (Model-View-Controller) ebenfalls ein Pattern, welches ich anwende. Es unterteilt die Anwendung in Models (Tabellen / Komponenten), Views (Ansichten, das GUI) und Controller, welche die CRUD Operationen enthalten. Oftmals verschiebt man die Logik in sogenannte Services. In meinem Projekt verwende ich einen AuthService, welche die Autorisierung/Authentifizierung \u00fcbernimmt. Ansonsten befindet sich die Logik bei mir in den Controllern.\",\n    Word.Style.normal\n);\n\n// Trailing empty paragraph.\naddParagraph(\"\", null);\n\n// Trailing empty, right-aligned paragraph.\nconst lastPara = addParagraph(\"\", null);\nlastPara.alignment = Word.Alignment.right;\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $app resolve to the Application, $d / $word.ActiveDocument to the document.\n#\n# Implements the changes from the commit \"Finished Project and added Privacy Page\":\n#  1. \"Erst beginne ich ...\" -> \"Erst danach beginne ich ...\" (real text change).\n#  2-5. Several paragraphs get their runs re-split around spell-checked words\n#       (Griptape, included, Product, Diagram, ...). The *visible* text of\n#       those paragraphs is byte-identical before/after (only <w:proofErr/>\n#       bookkeeping + run boundaries change, which Word inserts automatically\n#       during interactive spell-checking and is not reachable from the\n#       Word object model), so there is nothing observable to change there.\n#  6. A new \"Implementation\" section (Heading 1 + body paragraphs, a\n#     \"Problem mit Composite Pattern\" Heading 2 + body paragraphs, and two\n#     trailing empty paragraphs, the last right-aligned) is appended at the\n#     very end of the document body, right after the Class-Diagram picture\n#     and before the final section break.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) \"Zuerst habe ich mein Programm ... Erst beginne ich ...\"\n#    -> \"... Erst danach beginne ich ...\"\n# ---------------------------------------------------------------------\n$rng = $d.Content\n$rng.Find.Text = \"Erst beginne ich mit dem eigentlichen programmieren.\"\n$rng.Find.Replacement.Text = \"Erst danach beginne ich mit dem eigentlichen programmieren.\"\n$rng.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n# ---------------------------------------------------------------------\n# 6) Append the new \"Implementation\" section at the end of the document.\n# ---------------------------------------------------------------------\n\n# --- Paragraph: \"Implementation\" (Heading 1) ---\n$rng = $d.Content\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n$rng.Collapse(0)\n$rng.Text = \"Implementation\"\n$para = $d.Paragraphs.Last\n$para.Style = $d.Styles(\"berschrift1\")\n\n# --- Paragraph: intro sentence (Normal) ---\n$rng = $d.Content\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n$rng.Collapse(0)\n$rng.Text = \"Als Entwicklungsumgebung benutze ich Visual Studio 2019 mit einer MSSQL Datenbank. Ich erstelle ein ASP.NET Core MVC Projekt.\"\n$para = $d.Paragraphs.Last\n$para.Style = $d.Styles(\"Standard\")\n\n# --- Paragraph: \"Problem mit Composite Pattern\" (Heading 2) ---\n$rng = $d.Content\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n$rng.Collapse(0)\n$rng.Text = \"Problem mit Composite Pattern\"\n$para = $d.Paragraphs.Last\n$para.Style = $d.Styles(\"berschrift2\")\n\n# --- Paragraph: Composite Pattern problem description (Normal) ---\n$rng = $d.Content\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n$rng.Collapse(0)\n$rng.Text = \"Als ich das Composite Pattern implementieren wollte sah ich, dass es dabei ein Problem gibt. In einem MVC Projekt arbeitet man mit Model-Klassen, welche gleichzeitig die Tabellen in der Datenbank sind. Da eine Tabelle nicht von einer anderen Tabelle erben kann, war es nicht m\u00f6glich das Composite Pattern zu implementieren.\"\n$para = $d.Paragraphs.Last\n$para.Style = $d.Styles(\"Standard\")\n\n# --- Paragraph: Dependency Injection / MVC explanation (Normal) ---\n$rng = $d.Content\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n$rng.Collapse(0)\n$rng.Text = \"Ich habe beschlossen, stattdessen das Dependency Injection Design Pattern anzuwenden. Ausserdem ist MVC (Model-View-Controller) ebenfalls ein Pattern, welches ich anwende. Es unterteilt die Anwendung in Models (Tabellen / Komponenten), Views (Ansichten, das GUI) und Controller, welche die CRUD Operationen enthalten. Oftmals verschiebt man die Logik in sogenannte Services. In meinem Projekt verwende ich einen AuthService, welche die Autorisierung/Authentifizierung \u00fcbernimmt. Ansonsten befindet sich die Logik bei mir in den Controllern.\"\n$para = $d.Paragraphs.Last\n$para.Style = $d.Styles(\"Standard\")\n\n# --- Trailing empty paragraph ---\n$rng = $d.Content\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n$rng.Collapse(0)\n$para = $d.Paragraphs.Last\n$para.Style = $d.Styles(\"Standard\")\n\n# --- Trailing empty, right-aligned paragraph ---\n$rng = $d.Content\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n$rng.Collapse(0)\n$para = $d.Paragraphs.Last\n$para.Style = $d.Styles(\"Standard\")\n$para.Range.ParagraphFormat.Alignment = 2\n"}
